$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 113, shifting existing rows 113:243 down to 114:244.
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with its data.
$ws.Cells.Item(113, 1).Value = 7
$ws.Cells.Item(113, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(113, 3).Value = "Ñuble"
$ws.Cells.Item(113, 4).Value = [DateTime]"2022-08-08"
$ws.Cells.Item(113, 5).Value = 16
$ws.Cells.Item(113, 6).Value = 100112003
$ws.Cells.Item(113, 7).Value = "Ajo"
$ws.Cells.Item(113, 8).Value = "Chino"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 100
$ws.Cells.Item(113, 11).Value = 24000
$ws.Cells.Item(113, 12).Value = 25000
$ws.Cells.Item(113, 13).Value = 24500
$ws.Cells.Item(113, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(113, 15).Value = "China"
$ws.Cells.Item(113, 16).Value = 2450
$ws.Cells.Item(113, 17).Value = 10
$ws.Cells.Item(113, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(113, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
